$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 377.27585
$ws.Range("I33").Value = 343.73914
$ws.Range("J33").Value = 505.83334
$ws.Range("K33").Value = 343.73914
$ws.Range("L33").Value = 505.83334
$ws.Range("M33").Value = -114.73914
$ws.Range("N33").Value = -963.83334
$ws.Range("H62").Value = 8550477
$ws.Range("I62").Value = 11114609
$ws.Range("K62").Value = 11114609
$ws.Range("M62").Value = -11113985
$ws.Range("H65").Value = 8550477
$ws.Range("I65").Value = 11114609
$ws.Range("K65").Value = 55573045
$ws.Range("M65").Value = -55569925
$ws.Range("H74").Value = 3425
$ws.Range("I74").Value = 3580
$ws.Range("J74").Value = 3166.6667
$ws.Range("K74").Value = 3580
$ws.Range("L74").Value = 3166.6667
$ws.Range("M74").Value = -2644
$ws.Range("N74").Value = -5038.6667
$ws.Range("H77").Value = 3425
$ws.Range("I77").Value = 3580
$ws.Range("J77").Value = 3166.6667
$ws.Range("K77").Value = 17900
$ws.Range("L77").Value = 15833.3335
$ws.Range("M77").Value = -13220
$ws.Range("N77").Value = -25193.3335
$ws.Range("H129").Value = 712
$ws.Range("J129").Value = 831.625
$ws.Range("L129").Value = 2494.875
$ws.Range("N129").Value = -12494.875
$ws.Range("H132").Value = 9955.704
$ws.Range("I132").Value = 8052.087
$ws.Range("K132").Value = 24156.261
$ws.Range("M132").Value = -21626.261
$ws.Range("H137").Value = 1367.2106
$ws.Range("I137").Value = 888.34784
$ws.Range("J137").Value = 2101.4666
$ws.Range("K137").Value = 2665.04352
$ws.Range("L137").Value = 6304.399800000001
$ws.Range("M137").Value = -115.0435200000002
$ws.Range("N137").Value = -11404.3998

$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 5848.022
$ws.Range("I32").Value = 5755.778
$ws.Range("K32").Value = 5755.778
$ws.Range("M32").Value = -5468.778
$ws.Range("H45").Value = 1399.8
$ws.Range("I45").Value = 1249.75
$ws.Range("K45").Value = 1249.75
$ws.Range("M45").Value = -872.75
$ws.Range("H61").Value = 40001050
$ws.Range("I61").Value = 47619944
$ws.Range("K61").Value = 47619944
$ws.Range("M61").Value = -47619732
$ws.Range("H74").Value = 1056.4
$ws.Range("I74").Value = 1056.4
$ws.Range("K74").Value = 1056.4
$ws.Range("M74").Value = -182.4000000000001
$ws.Range("H77").Value = 1056.4
$ws.Range("I77").Value = 1056.4
$ws.Range("K77").Value = 5282
$ws.Range("M77").Value = -914
$ws.Range("H132").Value = 2711.926
$ws.Range("I132").Value = 2346.2354
$ws.Range("J132").Value = 3333.6
$ws.Range("K132").Value = 7038.706200000001
$ws.Range("L132").Value = 10000.8
$ws.Range("M132").Value = -4508.706200000001
$ws.Range("N132").Value = -15060.8
$ws.Range("H136").Value = 40001050
$ws.Range("I136").Value = 47619944
$ws.Range("K136").Value = 142859832
$ws.Range("M136").Value = -142857282

$ws = $wb.Worksheets.Item(3)
$ws.Range("H110").Value = 19175
$ws.Range("J110").Value = 19175
$ws.Range("L110").Value = 19175
$ws.Range("N110").Value = -27355
$ws.Range("H134").Value = 6554.1
$ws.Range("I134").Value = 1073.9375
$ws.Range("J134").Value = 28474.75
$ws.Range("K134").Value = 3221.8125
$ws.Range("L134").Value = 85424.25
$ws.Range("M134").Value = -686.8125
$ws.Range("N134").Value = -90494.25

$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 337.45456
$ws.Range("I7").Value = 172
$ws.Range("J7").Value = 432
$ws.Range("K7").Value = 172
$ws.Range("L7").Value = 432
$ws.Range("M7").Value = -59
$ws.Range("N7").Value = -658
$ws.Range("H31").Value = 1600.5834
$ws.Range("I31").Value = 1727.1111
$ws.Range("J31").Value = 1474.0555
$ws.Range("K31").Value = 1727.1111
$ws.Range("L31").Value = 1474.0555
$ws.Range("M31").Value = -1432.1111
$ws.Range("N31").Value = -2064.0555
$ws.Range("H34").Value = 1600.5834
$ws.Range("I34").Value = 1727.1111
$ws.Range("J34").Value = 1474.0555
$ws.Range("K34").Value = 1727.1111
$ws.Range("L34").Value = 1474.0555
$ws.Range("M34").Value = -1525.1111
$ws.Range("N34").Value = -1878.0555
$ws.Range("H62").Value = 3923917.2
$ws.Range("I62").Value = 2378.2979
$ws.Range("J62").Value = 50002000
$ws.Range("K62").Value = 2378.2979
$ws.Range("L62").Value = 50002000
$ws.Range("M62").Value = -1754.2979
$ws.Range("N62").Value = -50003248
$ws.Range("H65").Value = 3923917.2
$ws.Range("I65").Value = 2378.2979
$ws.Range("J65").Value = 50002000
$ws.Range("K65").Value = 11891.4895
$ws.Range("L65").Value = 250010000
$ws.Range("M65").Value = -8771.4895
$ws.Range("N65").Value = -250016240
$ws.Range("H132").Value = 2207.28
$ws.Range("I132").Value = 2033.4117
$ws.Range("K132").Value = 6100.2351
$ws.Range("M132").Value = -3570.2351

$ws = $wb.Worksheets.Item(5)
$ws.Range("H105").Value = 97228.09
$ws.Range("J105").Value = 97228.09
$ws.Range("L105").Value = 291684.27
$ws.Range("N105").Value = -296926.27
$ws.Range("H109").Value = 60018.35
$ws.Range("I109").Value = 72065.14
$ws.Range("K109").Value = 216195.42
$ws.Range("M109").Value = -215155.42
$ws.Range("H129").Value = 21930714
$ws.Range("I129").Value = 47619610
$ws.Range("J129").Value = 6945524
$ws.Range("K129").Value = 142858830
$ws.Range("L129").Value = 20836572
$ws.Range("M129").Value = -142853830
$ws.Range("N129").Value = -20846572

$ws = $wb.Worksheets.Item(6)
$ws.Range("H57").Value = 21500
$ws.Range("J57").Value = 21500
$ws.Range("L57").Value = 21500
$ws.Range("N57").Value = -23140
$ws.Range("H132").Value = 3131.318
$ws.Range("I132").Value = 3134.0715
$ws.Range("J132").Value = 3126.5
$ws.Range("K132").Value = 9402.2145
$ws.Range("L132").Value = 9379.5
$ws.Range("M132").Value = -6872.2145
$ws.Range("N132").Value = -14439.5

$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 606.1739
$ws.Range("I16").Value = 497.2857
$ws.Range("J16").Value = 1749.5
$ws.Range("K16").Value = 497.2857
$ws.Range("L16").Value = 1749.5
$ws.Range("M16").Value = -327.2857
$ws.Range("N16").Value = -2089.5
$ws.Range("H68").Value = 1998.75
$ws.Range("I68").Value = 1998.75
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1998.75
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1249.75
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 1998.75
$ws.Range("I71").Value = 1998.75
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9993.75
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -6249.75
$ws.Range("N71").ClearContents()
$ws.Range("H136").Value = 1071.5
$ws.Range("I136").Value = 950.5
$ws.Range("J136").Value = 1495
$ws.Range("K136").Value = 2851.5
$ws.Range("L136").Value = 4485
$ws.Range("M136").Value = -301.5
$ws.Range("N136").Value = -9585

$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 250010000
$ws.Range("I62").Value = 500000000
$ws.Range("J62").Value = 20003
$ws.Range("K62").Value = 500000000
$ws.Range("L62").Value = 20003
$ws.Range("M62").Value = -499999376
$ws.Range("N62").Value = -21251
$ws.Range("H65").Value = 250010000
$ws.Range("I65").Value = 500000000
$ws.Range("J65").Value = 20003
$ws.Range("K65").Value = 2500000000
$ws.Range("L65").Value = 100015
$ws.Range("M65").Value = -2499996880
$ws.Range("N65").Value = -106255
$ws.Range("H132").Value = 1764.44
$ws.Range("I132").Value = 1150.7894
$ws.Range("J132").Value = 3707.6667
$ws.Range("K132").Value = 3452.3682
$ws.Range("L132").Value = 11123.0001
$ws.Range("M132").Value = -922.3681999999999
$ws.Range("N132").Value = -16183.0001
$ws.Range("H136").Value = 1055.225
$ws.Range("I136").Value = 1010.36664
$ws.Range("K136").Value = 3031.09992
$ws.Range("M136").Value = -481.0999199999997
